# IP_calculation_Pro/equipment_database.xlsx
# "NEW Project IP_calculation_Pro add order number and count check"
#
# Sheet1 holds an equipment/IP-range table (columns: equipment_name,
# equipment_type, first_serial_number, first_IP, last_serial_number,
# last_IP). This edit:
#   - reworks the "Домик для Мышки Норушки" block (rows 3-5) so the
#     last_IP values are computed from a formula (256*19+255) instead of
#     hard-coded numbers, and resets the serial-number counters to 1
#   - fixes the last_IP on row 9 and resets the order/serial number on
#     row 10 to 1
#   - fixes last_IP on row 12 and resets the serial number on row 13 to 1
#   - moves the active selection to D16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: last_IP becomes a formula (256*19+255 = 5119)
$ws.Range("D3").Formula = "=256*19+255"

# Row 4: serial number reset to 1, last_IP same formula (shared range D4:D5)
$ws.Range("C4").Value = 1
$ws.Range("D4:D5").Formula = "=256*19+255"

# Row 5: serial number reset to 1 (last_IP already set via the D4:D5 fill above)
$ws.Range("C5").Value = 1

# Row 9: last_IP corrected to 5119
$ws.Range("D9").Value = 5119

# Row 10: order/serial number reset to 1
$ws.Range("C10").Value = 1

# Row 12: last_IP corrected to 511
$ws.Range("D12").Value = 511

# Row 13: serial number reset to 1
$ws.Range("C13").Value = 1

# Move the active selection to D16 (matches the saved view state)
$ws.Range("D16").Select()
